# Updates vm_pu.xlsx voltage-magnitude results for the "case with 380 kV" run.
# Columns B,C,D,E,F,I,J,K,L,M on rows 2-25 get new values; column G/H unaffected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("B","C","D","E","F","I","J","K","L","M")

$data = @(
    @(2, 1.02,1.072892327087843,1.07187021511377,1.076205767245857,1.084854449401716,1.052804418647453,1.077809583936984,1.074566238892569,1.078890311757843,1.087516390636249),
    @(3, 1.02,1.074397094034832,1.073022815664558,1.077513318269283,1.086160823543457,1.053211707738529,1.078970231807105,1.075534646404976,1.080014124295421,1.088640674898526),
    @(4, 1.02,1.075369739712431,1.073767481534872,1.078358574765836,1.087005272722473,1.053473160082334,1.079719760781944,1.076159544114202,1.080739936370977,1.089366743467912),
    @(5, 1.02,1.075778397506665,1.074080269215919,1.07871372982578,1.087360077275099,1.053582575919267,1.080034511184373,1.076421841486718,1.081044744174708,1.089671647303483),
    @(6, 1.02,1.075846998914589,1.07413277193459,1.078773350872027,1.087419638813767,1.053600918127456,1.080087338678587,1.076465858457818,1.081095903896159,1.089722822453484),
    @(7, 1.02,1.07537520116269,1.0737716620746,1.07836332110892,1.087010014425908,1.053474624059493,1.079723967867993,1.076163050549441,1.080744010493901,1.08937081892054),
    @(8, 1.02,1.073401089103901,1.072259981083573,1.076647832115718,1.085296125933076,1.0529424979584,1.078202141157856,1.074893876536371,1.079270395705259,1.087896644020756),
    @(9, 1.02,1.069914208633013,1.069587267643891,1.073618457584906,1.082269227434137,1.051988726531689,1.075508886204516,1.072644036389283,1.076663012365922,1.085287906042565),
    @(10, 1.02,1.067583680642252,1.067799196899214,1.071594231933962,1.080246429856688,1.051341938204248,1.073705290956814,1.071134899762848,1.074917297792348,1.083541048881349),
    @(11, 1.02,1.066573037837228,1.067023400473676,1.070716552022095,1.079369317412503,1.051059250038795,1.072922329509029,1.070479180912981,1.074159553250237,1.082782754423863),
    @(12, 1.02,1.066197405645003,1.066734997278332,1.070390359985739,1.079043328669163,1.050953850385978,1.072631197638986,1.07023527457366,1.07387781155888,1.082500800316375),
    @(13, 1.02,1.066277990760313,1.066796871586974,1.070460337643628,1.079113263056424,1.050976476958252,1.072693660313447,1.070287608916302,1.07393825893949,1.082561293627749),
    @(14, 1.02,1.066541992791012,1.066999565856873,1.070689592646249,1.079342374992804,1.051050545776585,1.07289827071562,1.070459026566696,1.074136270170666,1.08275945395145),
    @(15, 1.02,1.066704621854236,1.067124420817523,1.070830819843289,1.07948351305423,1.051096129384712,1.073024297333009,1.070564597024593,1.074258233944569,1.082881508509844),
    @(16, 1.02,1.067650721225467,1.067850650871159,1.071652455367563,1.080304614570398,1.051360643788427,1.073757211055134,1.071178369799199,1.074967547604568,1.083591334038891),
    @(17, 1.02,1.068243774941521,1.068305777306999,1.072167526444862,1.080819337162911,1.051525862349306,1.074216410905511,1.071562766806245,1.075411985274198,1.084036078116317),
    @(18, 1.02,1.068589548090895,1.068571095644893,1.072467845537891,1.081119448069274,1.051621978471062,1.074484062522349,1.071786762041964,1.075671041215486,1.084295307448207),
    @(19, 1.02,1.068707423430479,1.068661537182834,1.072570227570577,1.081221758276154,1.051654708731043,1.07457529243209,1.07186310195772,1.075759342731633,1.084383667177568),
    @(20, 1.02,1.068180160949973,1.068256962017677,1.072112275939056,1.080764124554131,1.051508162154744,1.074167162996764,1.071521547128487,1.075364319668844,1.083988380186829),
    @(21, 1.02,1.066464257307851,1.066939884052796,1.070622087850739,1.079274912532233,1.051028745325196,1.072838026521174,1.07040855785213,1.074077968616285,1.082701108717021),
    @(22, 1.02,1.065384040209875,1.066110406174395,1.069684087848483,1.07833748296549,1.050725020671367,1.072000577428419,1.069706789526602,1.073267556530596,1.081890070823553),
    @(23, 1.02,1.065956815287211,1.066550260451004,1.070181441791567,1.07883453855788,1.050886249312861,1.072444694568655,1.070079000196046,1.073697327828437,1.082320178294045),
    @(24, 1.02,1.068208905823033,1.06827902000708,1.072137241614375,1.080789073120635,1.051516160892884,1.074189416600089,1.071540173197483,1.075385858251474,1.084009933391832),
    @(25, 1.02,1.070816668567941,1.070279313222121,1.07440241994096,1.083052589558078,1.052237218986952,1.076206561943778,1.073227285605251,1.07733837793349,1.085963664588903)
)

foreach ($entry in $data) {
    $row = $entry[0]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $col = $cols[$i]
        $value = $entry[$i + 1]
        $ws.Range("$col$row").Value = $value
    }
}

$wb.Save()
